$wb = $excel.ActiveWorkbook

# Add the new "Backtracking" worksheet after the last existing sheet ("Heap")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Backtracking"

# Header row
$ws.Range("A1").Value = "Date Solved"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Algorithm"
$ws.Range("D1").Value = "Difficulty"
$ws.Range("E1").Value = "Solved First Time"
$ws.Range("F1").Value = "Video Help"
$ws.Range("G1").Value = "Revisit?"
$ws.Range("H1").Value = "Understand?"
$ws.Range("I1").Value = "Revisit Date #1"
$ws.Range("J1").Value = "Revisit Date #2"
$ws.Range("K1").Value = "Revisit Date #3"
$ws.Range("L1").Value = "Confidence Now"

# Data rows
# Format A2 as text first so the date-like string "10/24/25" is kept as a
# literal string instead of being auto-converted to a date serial number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "10/24/25"
$ws.Range("B2").Value = "Subsets"
$ws.Range("C2").Value = "Backtracking"
$ws.Range("D2").Value = "Medium"

$ws.Range("B3").Value = "Combination Sum"
$ws.Range("C3").Value = "Backtracking"
$ws.Range("D3").Value = "Medium"

$ws.Range("B4").Value = "Combination Sum 2"
$ws.Range("C4").Value = "Backtracking"
$ws.Range("D4").Value = "Medium"

$ws.Range("B5").Value = "Permutations"
$ws.Range("C5").Value = "Backtracking"
$ws.Range("D5").Value = "Medium"

$ws.Range("B6").Value = "Subsets 2"
$ws.Range("C6").Value = "Backtracking"
$ws.Range("D6").Value = "Medium"

$ws.Range("B7").Value = "Word Search"
$ws.Range("C7").Value = "Backtracking"
$ws.Range("D7").Value = "Medium"

$ws.Range("B8").Value = "Palindrome Partitioning"
$ws.Range("C8").Value = "Backtracking"
$ws.Range("D8").Value = "Medium"

$ws.Range("B9").Value = "Letter Combinations of a Phone Number"
$ws.Range("C9").Value = "Backtracking"
$ws.Range("D9").Value = "Medium"

$ws.Range("B10").Value = "N Queens"
$ws.Range("C10").Value = "Backtracking"
$ws.Range("D10").Value = "Hard"
